# Applies the "Updated cryptos list" price/volume refresh to Sheet1.
# Column D (Price) and E (Volume(1h)) are plain text cells (not numbers),
# so D-column updates are written with a leading apostrophe to keep Excel
# from re-typing strings like "208.10" as the number 208.1, then the cell
# style is reset to 'Normal' so no stray NumberFormat/quote-prefix style lingers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '26.622.27' },
    @{ Cell = 'E2'; Value = '  -0.29%  ' },
    @{ Cell = 'D3'; Value = '1.597.32' },
    @{ Cell = 'E3'; Value = '  -0.15%  ' },
    @{ Cell = 'E4'; Value = '  +0.26%  ' },
    @{ Cell = 'D5'; Value = '210.59' },
    @{ Cell = 'E7'; Value = '  +0.23%  ' },
    @{ Cell = 'D8'; Value = '0.0615' },
    @{ Cell = 'E8'; Value = '  -0.66%  ' },
    @{ Cell = 'D9'; Value = '0.246' },
    @{ Cell = 'E9'; Value = '  -0.66%  ' },
    @{ Cell = 'D10'; Value = '19.55' },
    @{ Cell = 'E10'; Value = '  +0.03%  ' },
    @{ Cell = 'D11'; Value = '0.0846' },
    @{ Cell = 'E11'; Value = '  +0.37%  ' },
    @{ Cell = 'D12'; Value = '1.820.85' },
    @{ Cell = 'E12'; Value = '  -0.20%  ' },
    @{ Cell = 'D13'; Value = '1.612.84' },
    @{ Cell = 'E13'; Value = '  +0.80%  ' },
    @{ Cell = 'E14'; Value = '  -0.03%  ' },
    @{ Cell = 'D15'; Value = '0.523' },
    @{ Cell = 'E15'; Value = '  -0.22%  ' },
    @{ Cell = 'D16'; Value = '64.64' },
    @{ Cell = 'E16'; Value = '  -0.99%  ' },
    @{ Cell = 'D17'; Value = '26.604.79' },
    @{ Cell = 'E17'; Value = '  -0.28%  ' },
    @{ Cell = 'E18'; Value = '  -2.33%  ' },
    @{ Cell = 'E19'; Value = '  +0.27%  ' },
    @{ Cell = 'D20'; Value = '208.10' },
    @{ Cell = 'E20'; Value = '  -0.87%  ' },
    @{ Cell = 'E21'; Value = '  -1.14%  ' },
    @{ Cell = 'E22'; Value = '  -0.20%  ' },
    @{ Cell = 'D23'; Value = '2.24' },
    @{ Cell = 'E23'; Value = '  -3.19%  ' },
    @{ Cell = 'D24'; Value = '8.93' },
    @{ Cell = 'E24'; Value = '  -0.09%  ' },
    @{ Cell = 'D25'; Value = '143.78' },
    @{ Cell = 'E25'; Value = '  +0.56%  ' },
    @{ Cell = 'E26'; Value = '  +0.23%  ' },
    @{ Cell = 'E27'; Value = '  +0.18%  ' },
    @{ Cell = 'E28'; Value = '  -1.02%  ' },
    @{ Cell = 'D29'; Value = '15.26' },
    @{ Cell = 'E29'; Value = '  -0.59%  ' },
    @{ Cell = 'D30'; Value = '0.0505' },
    @{ Cell = 'E30'; Value = '  -2.14%  ' },
    @{ Cell = 'E31'; Value = '  -0.37%  ' },
    @{ Cell = 'D32'; Value = '3.24' },
    @{ Cell = 'E32'; Value = '  -0.48%  ' },
    @{ Cell = 'E33'; Value = '  -0.62%  ' },
    @{ Cell = 'E34'; Value = '  +19.27%  ' },
    @{ Cell = 'D35'; Value = '1.277.30' },
    @{ Cell = 'E35'; Value = '  -1.07%  ' },
    @{ Cell = 'E36'; Value = '  +0.90%  ' },
    @{ Cell = 'D37'; Value = '1.49' },
    @{ Cell = 'E37'; Value = '  -0.65%  ' },
    @{ Cell = 'D38'; Value = '0.597' },
    @{ Cell = 'E38'; Value = '  -3.62%  ' },
    @{ Cell = 'E39'; Value = '  -2.17%  ' },
    @{ Cell = 'D40'; Value = '0.821' },
    @{ Cell = 'E40'; Value = '  -0.46%  ' },
    @{ Cell = 'E41'; Value = '  +0.04%  ' },
    @{ Cell = 'B42'; Value = 'MXToken' },
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' },
    @{ Cell = 'D42'; Value = '2.16' },
    @{ Cell = 'E42'; Value = '  -1.03%  ' },
    @{ Cell = 'B43'; Value = 'TrustWalletToken' },
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' },
    @{ Cell = 'D43'; Value = '0.773' },
    @{ Cell = 'E43'; Value = '  -1.44%  ' },
    @{ Cell = 'D44'; Value = '62.53' },
    @{ Cell = 'E44'; Value = '  -1.04%  ' },
    @{ Cell = 'D45'; Value = '1.733.09' },
    @{ Cell = 'E45'; Value = '  -0.35%  ' },
    @{ Cell = 'D46'; Value = '89.62' },
    @{ Cell = 'E46'; Value = '  -1.52%  ' },
    @{ Cell = 'D47'; Value = '1.57' },
    @{ Cell = 'E47'; Value = '  -0.52%  ' },
    @{ Cell = 'D48'; Value = '0.0₆0105' },
    @{ Cell = 'E48'; Value = '  -1.64%  ' },
    @{ Cell = 'E49'; Value = '  +1.89%  ' },
    @{ Cell = 'D50'; Value = '0.0512' },
    @{ Cell = 'E50'; Value = '  +0.59%  ' },
    @{ Cell = 'E51'; Value = '  +0.13%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Cell.StartsWith('D')) {
        $cell.Value = "'" + $u.Value
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $u.Value
    }
}
